# Update "horarios" workbook with newest scrape batch (timestamp 07:13:03)
# Adds new rows to each of the three sheets and updates the header summary
# cells ("Ultima actualizacion" / "Total filas") accordingly.

$wb = $excel.ActiveWorkbook

$newTimestamp = "Última actualización: 07:13:03"
$scrapTime = "07:13:03"

# ---------------------------------------------------------------------------
# Sheet 1: LP1912  (rows 6..39 existing -> add rows 40..42)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = $newTimestamp
$ws1.Range("A3").Value = "Total filas: 37"

$ws1.Cells.Item(40, 1).Value = $scrapTime
$ws1.Cells.Item(40, 2).Value = "08:52"
$ws1.Cells.Item(40, 3).Value = "23_HERNANDEZ"
$ws1.Cells.Item(40, 4).Value = 99
$ws1.Cells.Item(40, 5).Value = "LP1912"

$ws1.Cells.Item(41, 1).Value = $scrapTime
$ws1.Cells.Item(41, 2).Value = "08:54"
$ws1.Cells.Item(41, 3).Value = "215B_EL PATO"
$ws1.Cells.Item(41, 4).Value = 101
$ws1.Cells.Item(41, 5).Value = "LP1912"

$ws1.Cells.Item(42, 1).Value = $scrapTime
$ws1.Cells.Item(42, 2).Value = "08:58"
$ws1.Cells.Item(42, 3).Value = "215A_EL PATO"
$ws1.Cells.Item(42, 4).Value = 105
$ws1.Cells.Item(42, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215  (rows 6..14 existing -> add rows 15..16)
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = $newTimestamp
$ws2.Range("A3").Value = "Total filas: 11"

$ws2.Cells.Item(15, 1).Value = $scrapTime
$ws2.Cells.Item(15, 2).Value = "08:54"
$ws2.Cells.Item(15, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(15, 4).Value = 101
$ws2.Cells.Item(15, 5).Value = "LP1912"

$ws2.Cells.Item(16, 1).Value = $scrapTime
$ws2.Cells.Item(16, 2).Value = "08:58"
$ws2.Cells.Item(16, 3).Value = "215A_EL PATO"
$ws2.Cells.Item(16, 4).Value = 105
$ws2.Cells.Item(16, 5).Value = "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173  (rows 6..10 existing -> add row 11)
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = $newTimestamp
$ws3.Range("A3").Value = "Total filas: 6"

$ws3.Cells.Item(11, 1).Value = $scrapTime
$ws3.Cells.Item(11, 2).Value = "08:51"
$ws3.Cells.Item(11, 3).Value = "215C_LA PLATA"
$ws3.Cells.Item(11, 4).Value = 98
$ws3.Cells.Item(11, 5).Value = "L6203"
